$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new first column (A) holding a numeric row id, shifting the
# previous header texts one column to the right, and reuse the
# "NAMA TAMU" style text ("dfdfldsfhl") in the newly vacated column F.
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "2024-10-07"
$ws.Range("C1").Value = "szads"
$ws.Range("D1").Value = "dfdfldsfhl"
$ws.Range("E1").Value = "0987654"
$ws.Range("F1").Value = "dfdfldsfhl"
$ws.Range("G1").Value = "s"

$wb.Save()
